# Update the "想去人数" (attendance count) figures that changed between
# data refreshes on the gh-pages generated output.
#
# Sheet "展览"   (sheet1): F2 13630->13643, F3 321->322, F6 476->477, F7 1383->1389
# Sheet "全部类型" (sheet4): F2 13630->13643, F3 321->322, F8 476->477, F9 1383->1389

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13643
$ws1.Range("F3").Value = 322
$ws1.Range("F6").Value = 477
$ws1.Range("F7").Value = 1389

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13643
$ws4.Range("F3").Value = 322
$ws4.Range("F8").Value = 477
$ws4.Range("F9").Value = 1389
